$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = "'3"
$ws.Range("L2").Value = "maa://39402 (94.87), *maa://30515 (70.37), *maa://34787 (73.91)"
$ws.Range("S2").Value = "'1"
$ws.Range("T2").Value = "maa://22742 (91.47)"
$ws.Range("AA2").Value = "'2"
$ws.Range("AB2").Value = "maa://21246 (91.48), maa://36684 (93.12)"
$ws.Range("AE2").Value = "'2"
$ws.Range("AF2").Value = "maa://25251 (91.91), maa://59087 (100.0)"
$ws.Range("G3").Value = "'1"
$ws.Range("H3").Value = "maa://21247 (98.44)"
$ws.Range("O3").Value = "'2"
$ws.Range("P3").Value = "maa://21249 (94.68), maa://26254 (97.22)"
$ws.Range("S3").Value = "'2"
$ws.Range("T3").Value = "maa://24617 (90.91), maa://45854 (84.78)"
$ws.Range("C4").Value = "'3"
$ws.Range("D4").Value = "maa://24632 (94.09), maa://22499 (86.67), maa://22746 (100.0)"
$ws.Range("S4").Value = "'4"
$ws.Range("T4").Value = "maa://32509 (93.98), maa://27295 (88.3), maa://22754 (89.19), *maa://31008 (79.55)"
$ws.Range("W4").Value = "'1"
$ws.Range("X4").Value = "maa://43217 (93.8)"
$ws.Range("AE7").Value = "'1"
$ws.Range("AF7").Value = "maa://45272 (97.06)"
$ws.Range("A8").Value = "更新日期：2025.06.10 14:12:09"
$ws.Range("O8").Value = "'3"
$ws.Range("P8").Value = "maa://32931 (83.8), maa://23252 (91.67), maa://37496 (97.44)"
$ws.Range("G9").Value = "'1"
$ws.Range("H9").Value = "maa://56348 (100.0)"
$ws.Range("S9").Value = "'1"
$ws.Range("T9").Value = "maa://26222 (98.33)"
$ws.Range("AA9").Value = "'2"
$ws.Range("AB9").Value = "maa://28711 (87.23), maa://40166 (93.75)"
$ws.Range("S10").Value = "'2"
$ws.Range("T10").Value = "maa://27395 (96.73), maa://22755 (89.15)"
$ws.Range("T11").Value = "maa://22747 (90.96), maa://22501 (98.23), maa://45521 (91.67)"
$ws.Range("AA11").Value = "'2"
$ws.Range("AB11").Value = "maa://29912 (97.78), maa://22516 (87.36)"
$ws.Range("G12").Value = "'2"
$ws.Range("H12").Value = "maa://21867 (90.81), maa://54294 (100.0)"
$ws.Range("O13").Value = "'3"
$ws.Range("P13").Value = "maa://22676 (93.66), *maa://22583 (78.05), maa://48321 (92.31)"
$ws.Range("O14").Value = "'3"
$ws.Range("P14").Value = "maa://23250 (98.86), maa://20107 (87.1), maa://22772 (100.0)"
$ws.Range("AE16").Value = "'1"
$ws.Range("AF16").Value = "maa://27755 (94.06)"
$ws.Range("H17").Value = "maa://22430 (89.67), maa://39599 (86.42)"
$ws.Range("O17").Value = "'2"
$ws.Range("P17").Value = "maa://23890 (81.08), maa://56238 (100.0)"
$ws.Range("K18").Value = "'2"
$ws.Range("L18").Value = "maa://22466 (92.12), maa://52226 (95.65)"
$ws.Range("C20").Value = "'3"
$ws.Range("D20").Value = "maa://21432 (90.73), maa://25198 (93.89), maa://36680 (91.18)"
$ws.Range("G20").Value = "'1"
$ws.Range("H20").Value = "maa://22864 (90.81)"
$ws.Range("G22").Value = "'1"
$ws.Range("H22").Value = "maa://25236 (96.19)"
$ws.Range("L23").Value = "maa://39756 (95.77), maa://39875 (94.81)"
$ws.Range("O23").Value = "'3"
$ws.Range("P23").Value = "maa://30587 (92.09), *maa://29748 (76.3), *maa://37566 (78.26)"
$ws.Range("W24").Value = "'5"
$ws.Range("X24").Value = "maa://29988 (85.05), maa://23504 (93.67), *maa://25141 (77.37), *maa://36663 (78.0), maa://52227 (100.0)"
$ws.Range("AE24").Value = "'4"
$ws.Range("AF24").Value = "maa://22523 (81.74), *maa://36672 (75.38), maa://29910 (93.75), maa://45831 (85.71)"
$ws.Range("H25").Value = "*maa://29063 (72.82), *maa://25311 (74.11), ***maa://22725 (4.76), *maa://45047 (73.33)"
$ws.Range("G26").Value = "'1"
$ws.Range("H26").Value = "maa://24913 (91.26)"
$ws.Range("O26").Value = "'2"
$ws.Range("P26").Value = "maa://39870 (92.86), maa://56625 (100.0)"
$ws.Range("G27").Value = "'2"
$ws.Range("H27").Value = "*maa://39601 (79.17), maa://34494 (97.22)"
$ws.Range("D28").Value = "maa://24465 (90.65), maa://25725 (82.83)"
$ws.Range("L28").Value = "maa://30770 (82.0)"
$ws.Range("S28").Value = "'1"
$ws.Range("T28").Value = "maa://23263 (95.45)"
$ws.Range("W28").Value = "'2"
$ws.Range("X28").Value = "maa://39929 (91.95), maa://41749 (92.25)"
$ws.Range("AE28").Value = "'1"
$ws.Range("AF28").Value = "maa://36660 (92.54)"
$ws.Range("O29").Value = "'1"
$ws.Range("P29").Value = "maa://54169 (100.0)"
$ws.Range("G31").Value = "'1"
$ws.Range("H31").Value = "maa://32721 (100.0)"
$ws.Range("L31").Value = "maa://35926 (93.73), maa://36258 (87.8), *maa://43904 (73.33)"
$ws.Range("G32").Value = "'3"
$ws.Range("H32").Value = "maa://21895 (97.62), maa://36667 (97.32), maa://22760 (100.0)"
$ws.Range("P34").Value = "maa://48817 (96.83), maa://56235 (100.0)"
$ws.Range("L35").Value = "maa://41296 (97.35)"
$ws.Range("L37").Value = "maa://45718 (98.44), maa://47069 (81.82), maa://56336 (87.5), maa://45789 (100.0)"
$ws.Range("G39").Value = "'4"
$ws.Range("H39").Value = "maa://36670 (89.57), maa://25199 (85.22), maa://30434 (93.13), *maa://45059 (79.41)"
$ws.Range("G45").Value = "'3"
$ws.Range("H45").Value = "maa://21229 (84.42), maa://30807 (93.51), maa://42459 (89.47)"
$ws.Range("G53").Value = "'1"
$ws.Range("H53").Value = "maa://32534 (95.12)"
